$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append to the data table (CONTAD, CONFIRMADOS_UTI, DESCARTADOS_UTI,
# AGUARDANDO_UTI, LEITOSDISP_UTI, PACIENTES_TOTAL, LEITOSTOAL_UTI)
$newRows = @(
    @(46, 62, 2, 12, 29, 76, 105),
    @(47, 62, 2, 12, 29, 76, 105),
    @(48, 75, 3, 12, 20, 90, 110)
)

$startRow = 47
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowValues = $newRows[$i]
    for ($c = 1; $c -le $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
